$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$updates = @(
    "M70=-5728.5",
    "I70=1999.5",
    "L70=16580.538",
    "H70=5056.533",
    "K70=5998.5",
    "J70=5526.846",
    "N70=-17120.538",
    "K73=5998.5",
    "M73=-5062.5",
    "H73=5056.533",
    "N73=-18452.538",
    "L73=16580.538",
    "J73=5526.846",
    "I73=1999.5",
    "H132=2013.75",
    "I132=2114.6667",
    "K132=6344.000100000001",
    "M132=-3814.000100000001",
    "N137=-18421.8",
    "K137=6196.5",
    "I137=2065.5",
    "L137=13321.8",
    "H137=2857.2",
    "J137=4440.6",
    "M137=-3646.5",
    "M138=1967.8462",
    "N138=-25521.44",
    "L138=15241.44",
    "J138=5080.48",
    "K138=3172.1538",
    "H138=3704.158",
    "I138=1057.3846"
)
foreach ($u in $updates) {
    $parts = $u.Split("=")
    $ws.Range($parts[0]).Value = [double]$parts[1]
}

$ws = $wb.Worksheets.Item("ARM")
$updates = @(
    "M2=-4886",
    "H2=4999.3335",
    "K2=4999",
    "I2=4999",
    "I32=4990",
    "H32=7242.25",
    "K32=4990",
    "M32=-4703",
    "I116=4999",
    "M116=-2705",
    "H116=4999.3335",
    "K116=4999",
    "L122=8614.667099999999",
    "N122=-13514.6671",
    "M122=-1584891.25",
    "H122=359964.47",
    "I122=529113.75",
    "J122=2871.5557",
    "K122=1587341.25",
    "L132=6748.5",
    "H132=2249.5",
    "J132=2249.5",
    "I132=2249.5",
    "K132=6748.5",
    "N132=-11808.5",
    "M132=-4218.5"
)
foreach ($u in $updates) {
    $parts = $u.Split("=")
    $ws.Range($parts[0]).Value = [double]$parts[1]
}

$ws = $wb.Worksheets.Item("BSM")
$updates = @(
    "K3=4999",
    "M3=-4885",
    "I3=4999",
    "H3=4999.3335",
    "I68=5000",
    "L68=0",
    "H68=5000",
    "J68=0",
    "M68=-4189",
    "K68=5000",
    "I71=5000",
    "L71=0",
    "H71=5000",
    "M71=-10944",
    "J71=0",
    "K71=15000",
    "M86=-2320",
    "H86=4564.3335",
    "I86=3443",
    "K86=3443",
    "H89=4564.3335",
    "M89=-11599",
    "I89=3443",
    "K89=17215",
    "N134=-16455",
    "L134=11385",
    "K134=7809",
    "I134=2603",
    "J134=3795",
    "M134=-5274",
    "H134=3099.6667"
)
foreach ($u in $updates) {
    $parts = $u.Split("=")
    $ws.Range($parts[0]).Value = [double]$parts[1]
}
$ws.Range("N68").ClearContents()
$ws.Range("N71").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$updates = @(
    "L16=2000",
    "H16=1455.5",
    "J16=2000",
    "N16=-2574",
    "L31=4319.609",
    "H31=4088.55",
    "J31=4319.609",
    "N31=-4909.609",
    "H34=4088.55",
    "J34=4319.609",
    "N34=-4723.609",
    "L34=4319.609",
    "M105=747.6",
    "K105=999.4",
    "I105=999.4",
    "H105=1288.1428",
    "J113=2000",
    "L113=2000",
    "H113=1455.5",
    "N113=-6340",
    "L132=29994",
    "H132=2042.9778",
    "J132=9998",
    "I132=1672.9767",
    "K132=5018.9301",
    "N132=-35054",
    "M132=-2488.9301",
    "N134=-16904.625",
    "L134=11834.625",
    "K134=5577.9231",
    "I134=1859.3077",
    "J134=3944.875",
    "M134=-3042.9231",
    "H134=2350.0293"
)
foreach ($u in $updates) {
    $parts = $u.Split("=")
    $ws.Range($parts[0]).Value = [double]$parts[1]
}

$ws = $wb.Worksheets.Item("CUL")
$updates = @(
    "N4=-75566.375",
    "M4=-279345.65",
    "J4=25114.125",
    "K4=279457.65",
    "I4=93152.55",
    "L4=75342.375",
    "H4=79195.94500000001",
    "H131=1478.9"
)
foreach ($u in $updates) {
    $parts = $u.Split("=")
    $ws.Range($parts[0]).Value = [double]$parts[1]
}

$ws = $wb.Worksheets.Item("GSM")
$updates = @(
    "L122=609411.6000000001",
    "N122=-614311.6000000001",
    "M122=-5660.3638",
    "H122=65339",
    "I122=2703.4546",
    "J122=203137.2",
    "K122=8110.3638",
    "L132=8737.071599999999",
    "H132=2502.423",
    "J132=2912.3572",
    "I132=2024.1666",
    "K132=6072.4998",
    "N132=-13797.0716",
    "M132=-3542.4998"
)
foreach ($u in $updates) {
    $parts = $u.Split("=")
    $ws.Range($parts[0]).Value = [double]$parts[1]
}

$ws = $wb.Worksheets.Item("LTW")
$updates = @(
    "N42=-19751126",
    "L42=19750000",
    "J42=19750000",
    "H42=19750000",
    "L46=0",
    "I46=0",
    "J46=0",
    "K46=0",
    "H46=0",
    "J49=19750000",
    "N49=-19750294",
    "L49=19750000",
    "H49=19750000",
    "J136=5894.5",
    "H136=4144.8945",
    "L136=17683.5",
    "N136=-22783.5"
)
foreach ($u in $updates) {
    $parts = $u.Split("=")
    $ws.Range($parts[0]).Value = [double]$parts[1]
}
$ws.Range("N46").ClearContents()
$ws.Range("M46").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$updates = @(
    "I62=8049.5",
    "M62=-7425.5",
    "K62=8049.5",
    "H62=8013.8125",
    "H65=8013.8125",
    "M65=-37127.5",
    "I65=8049.5",
    "K65=40247.5",
    "H132=87110.09",
    "I132=135788.58",
    "K132=407365.74",
    "M132=-404835.74"
)
foreach ($u in $updates) {
    $parts = $u.Split("=")
    $ws.Range($parts[0]).Value = [double]$parts[1]
}
